# Segundo intento con menos variables
#
# Adds a new worksheet "sinhuevo" (per-item objective/cofre breakdown for
# "pan"/"galleta", same visual language as the existing
# "speedcraft-alm-cofre" sheet) as the last tab, and refreshes the column
# widths / selection on "speedcraft-alm-cofre" itself.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. speedcraft-alm-cofre: add a narrow spacer column A and re-fit
#    columns C:F, then leave the selection on whole columns A:XFD
#    (mirrors a "select all columns" action) instead of A5:K5.
# ---------------------------------------------------------------------
$alm = $wb.Worksheets.Item("speedcraft-alm-cofre")

$alm.Columns.Item(1).ColumnWidth = 2
$alm.Columns.Item(3).ColumnWidth = 13.5
$alm.Columns.Item(4).ColumnWidth = 14.166666666666666
$alm.Columns.Item(5).ColumnWidth = 13.5
$alm.Columns.Item(6).ColumnWidth = 12.666666666666666

$alm.Range("A1:XFD1048576").Select()

# ---------------------------------------------------------------------
# 2. Add the new sheet "sinhuevo" as the last tab (after
#    speedcraft-alm-cofre) and populate it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet, $null, $null)
$new.Name = "sinhuevo"

# -- headers (row 1) --
$new.Cells.Item(1,1).Value = "a-objetivo"
$new.Cells.Item(1,2).Value = "alm1-cacao"
$new.Cells.Item(1,3).Value = "alm1-leche"
$new.Cells.Item(1,4).Value = "alm1-trigo"
$new.Cells.Item(1,5).Value = "cofre-galleta"
$new.Cells.Item(1,6).Value = "cofre-galleta-cacao"
$new.Cells.Item(1,7).Value = "cofre-galleta-leche"
$new.Cells.Item(1,8).Value = "cofre-galleta-trigo"
$new.Cells.Item(1,9).Value = "cofre-pan"
$new.Cells.Item(1,10).Value = "cofre-pan-cacao"
$new.Cells.Item(1,11).Value = "cofre-pan-leche"
$new.Cells.Item(1,12).Value = "cofre-pan-trigo"

# -- data rows 2-6 --
$new.Cells.Item(2,1).Value = "pan"
$new.Cells.Item(2,2).Value = 15
$new.Cells.Item(2,3).Value = 20
$new.Cells.Item(2,4).Value = 30
$new.Cells.Item(2,5).Value = 0
$new.Cells.Item(2,6).Value = 0
$new.Cells.Item(2,7).Value = 0
$new.Cells.Item(2,8).Value = 0
$new.Cells.Item(2,9).Value = 0
$new.Cells.Item(2,10).Value = 0
$new.Cells.Item(2,11).Value = 0
$new.Cells.Item(2,12).Value = 0

$new.Cells.Item(3,1).Value = "pan"
$new.Cells.Item(3,2).Value = 15
$new.Cells.Item(3,3).Value = 20
$new.Cells.Item(3,4).Value = 20
$new.Cells.Item(3,5).Value = 0
$new.Cells.Item(3,6).Value = 0
$new.Cells.Item(3,7).Value = 0
$new.Cells.Item(3,8).Value = 0
$new.Cells.Item(3,9).Value = 0
$new.Cells.Item(3,10).Value = 0
$new.Cells.Item(3,11).Value = 0
$new.Cells.Item(3,12).Value = 10

$new.Cells.Item(4,1).Value = "galleta"
$new.Cells.Item(4,2).Value = 15
$new.Cells.Item(4,3).Value = 15
$new.Cells.Item(4,4).Value = 20
$new.Cells.Item(4,5).Value = 0
$new.Cells.Item(4,6).Value = 0
$new.Cells.Item(4,7).Value = 0
$new.Cells.Item(4,8).Value = 0
$new.Cells.Item(4,9).Value = 5
$new.Cells.Item(4,10).Value = 0
$new.Cells.Item(4,11).Value = 0
$new.Cells.Item(4,12).Value = 0

$new.Cells.Item(5,1).Value = "galleta"
$new.Cells.Item(5,2).Value = 5
$new.Cells.Item(5,3).Value = 15
$new.Cells.Item(5,4).Value = 20
$new.Cells.Item(5,5).Value = 0
$new.Cells.Item(5,6).Value = 5
$new.Cells.Item(5,7).Value = 0
$new.Cells.Item(5,8).Value = 0
$new.Cells.Item(5,9).Value = 5
$new.Cells.Item(5,10).Value = 0
$new.Cells.Item(5,11).Value = 0
$new.Cells.Item(5,12).Value = 0

$new.Cells.Item(6,1).Value = "galleta"
$new.Cells.Item(6,2).Value = 5
$new.Cells.Item(6,3).Value = 15
$new.Cells.Item(6,4).Value = 10
$new.Cells.Item(6,5).Value = 5
$new.Cells.Item(6,6).Value = 0
$new.Cells.Item(6,7).Value = 0
$new.Cells.Item(6,8).Value = 0
$new.Cells.Item(6,9).Value = 5
$new.Cells.Item(6,10).Value = 0
$new.Cells.Item(6,11).Value = 0
$new.Cells.Item(6,12).Value = 0

# ---------------------------------------------------------------------
# 3. Formats: copy the existing block-shading/border styles from
#    speedcraft-alm-cofre (same visual language) onto the matching
#    columns/rows of the new sheet, reusing the existing style records.
# ---------------------------------------------------------------------
$alm.Range("C1").Copy()
$new.Range("B1:D1").PasteSpecial(-4122)

$alm.Range("G1").Copy()
$new.Range("E1:H1").PasteSpecial(-4122)

$alm.Range("J8").Copy()
$new.Range("I1:L1").PasteSpecial(-4122)

$alm.Range("C5").Copy()
$new.Range("B2:D2").PasteSpecial(-4122)

$alm.Range("G5").Copy()
$new.Range("E2:H2").PasteSpecial(-4122)

$alm.Range("J7").Copy()
$new.Range("I2:L2").PasteSpecial(-4122)

$alm.Range("C1").Copy()
$new.Range("B3:D3").PasteSpecial(-4122)

$alm.Range("G1").Copy()
$new.Range("E3:H3").PasteSpecial(-4122)

$alm.Range("J8").Copy()
$new.Range("I3:L3").PasteSpecial(-4122)

$alm.Range("C5").Copy()
$new.Range("B4:D4").PasteSpecial(-4122)

$alm.Range("G5").Copy()
$new.Range("E4:H4").PasteSpecial(-4122)

$alm.Range("J7").Copy()
$new.Range("I4:L4").PasteSpecial(-4122)

$alm.Range("B5").Copy()
$new.Range("A5").PasteSpecial(-4122)

$alm.Range("C1").Copy()
$new.Range("B5:D5").PasteSpecial(-4122)

$alm.Range("G1").Copy()
$new.Range("E5:H5").PasteSpecial(-4122)

$alm.Range("J8").Copy()
$new.Range("I5:L5").PasteSpecial(-4122)

$alm.Range("C1").Copy()
$new.Range("B6:D6").PasteSpecial(-4122)

$alm.Range("G1").Copy()
$new.Range("E6:H6").PasteSpecial(-4122)

$alm.Range("J8").Copy()
$new.Range("I6:L6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Column widths for the new sheet (nearest achievable to the source
#    workbook's own bestFit measurements).
# ---------------------------------------------------------------------
$new.Columns.Item(1).ColumnWidth = 5.666666666666667
$new.Columns.Item(2).ColumnWidth = 10
$new.Columns.Item(3).ColumnWidth = 10
$new.Columns.Item(4).ColumnWidth = 9.333333333333334
$new.Columns.Item(5).ColumnWidth = 11.5
$new.Columns.Item(6).ColumnWidth = 17.166666666666668
$new.Columns.Item(7).ColumnWidth = 17.166666666666668
$new.Columns.Item(8).ColumnWidth = 16.5
$new.Columns.Item(9).ColumnWidth = 8.666666666666666
$new.Columns.Item(10).ColumnWidth = 14.333333333333334
$new.Columns.Item(11).ColumnWidth = 14.333333333333334
$new.Columns.Item(12).ColumnWidth = 13.666666666666666

# ---------------------------------------------------------------------
# 5. Selection + active sheet.
# ---------------------------------------------------------------------
$new.Range("B5:L5").Select()
$new.Activate()

Write-Output "edit complete"
